# Auto update Excel log
# Appends new sensor-log rows to the mmWave, PIR and Humidity sheets,
# mirroring the automated logger that produced this workbook.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $Sheet,
        [int]$Row,
        [string]$Date,
        [string]$Time,
        [string]$Hour,
        [string]$Location,
        [string]$Value,
        [string]$Status
    )

    # Columns A (date-looking) and E (can look like a percentage) need to be
    # forced to text so Excel doesn't silently convert them to a date serial
    # / percentage number - the source log stores everything as plain text.
    $Sheet.Cells.Item($Row, 1).NumberFormat = "@"
    $Sheet.Cells.Item($Row, 1).Value = $Date
    $Sheet.Cells.Item($Row, 2).Value = $Time
    $Sheet.Cells.Item($Row, 3).Value = $Hour
    $Sheet.Cells.Item($Row, 4).Value = $Location
    $Sheet.Cells.Item($Row, 5).NumberFormat = "@"
    $Sheet.Cells.Item($Row, 5).Value = $Value
    $Sheet.Cells.Item($Row, 6).Value = $Status
}

# ---- mmWave sheet: append rows 12-22 ----
$mmWave = $wb.Worksheets.Item("mmWave")

Add-LogRow $mmWave 12 "2026-01-30" "14:50:57" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmWave 13 "2026-01-30" "14:51:08" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmWave 14 "2026-01-30" "14:51:18" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmWave 15 "2026-01-30" "14:51:29" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmWave 16 "2026-01-30" "14:52:10" "14:00" "Living Room" "CRITICAL EMERGENCY" "FALL_DETECTED"
Add-LogRow $mmWave 17 "2026-01-30" "14:52:14" "14:00" "Living Room" "CRITICAL EMERGENCY" "FALL_DETECTED"
Add-LogRow $mmWave 18 "2026-01-30" "14:52:14" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmWave 19 "2026-01-30" "14:52:18" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmWave 20 "2026-01-30" "14:52:28" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmWave 21 "2026-01-30" "14:52:45" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $mmWave 22 "2026-01-30" "14:53:00" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"

# ---- PIR sheet: append rows 4-6 ----
$pir = $wb.Worksheets.Item("PIR")

Add-LogRow $pir 4 "2026-01-30" "14:50:55" "14:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 5 "2026-01-30" "14:50:56" "14:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 6 "2026-01-30" "14:51:01" "14:00" "Bathroom" "No Motion" "Inactive"

# ---- Humidity sheet: append rows 4-5 ----
$humidity = $wb.Worksheets.Item("Humidity")

Add-LogRow $humidity 4 "2026-01-30" "14:50:55" "14:00" "Bathroom" "88.3%" "Active"
Add-LogRow $humidity 5 "2026-01-30" "14:51:02" "14:00" "Bathroom" "88.7%" "Active"
